$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.890.55'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '2.230.92'
$ws.Range("E3").Value = '  -4.50%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.52'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '84.82'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  -2.80%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -3.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0783'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.82'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.15%  '
$ws.Range("E12").Value = '  -10.82%  '
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").Value = '2.576.19'
$ws.Range("E14").Value = '  -4.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.32'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.17'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.34%  '
$ws.Range("D17").Value = '2.228.94'
$ws.Range("E17").Value = '  -5.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.721'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.85%  '
$ws.Range("D19").Value = '39.804.52'
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D20").Value = '0.0₃0879'
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("E21").Value = '  -5.48%  '
$ws.Range("E22").Value = '  -3.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.46'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '234.66'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.47%  '
$ws.Range("E27").Value = '  +0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.83'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.22'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.53'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '150.02'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.85'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.22%  '
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0703'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.02%  '
$ws.Range("E37").Value = '  +5.90%  '
$ws.Range("E38").Value = '  -2.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0979'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.82%  '
$ws.Range("E40").Value = '  -2.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.67'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.68'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.27%  '
$ws.Range("D43").Value = '1.933.35'
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("E44").Value = '  -2.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0266'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.24'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.51'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.62'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.37%  '
$ws.Range("D49").Value = '2.446.10'
$ws.Range("E49").Value = '  -4.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.27'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.80'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.30%  '
